$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

# New match-day rows (fecha 2025-08-09 / serial 45878) appended to the Partidos log.
$data = @(
    @(45878, "Fabian Caicedo",            "Amarillo", "Arquero",       0, 0, $true,  1, 0, 0, 0, 0),
    @(45878, "Invitado",                  "Azul",     "Arquero",       0, 0, $true,  5, 0, 0, 0, 0),
    @(45878, "Carlos Fernando Valencia",  "Amarillo", "Delantero",     4, 0, $false, 0, 0, 0, 0, 0),
    @(45878, "Armando Murillo",           "Amarillo", "Defensa",       1, 0, $false, 0, 0, 0, 1, 0),
    @(45878, "Juan Diego Gomez Ceballos", "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @(45878, "Harold Gonzalez Castro",    "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(45878, "Juan David Espinal",        "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @(45878, "Alexander Uribe",           "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0)
)

$startRow = 373
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = $row[5]
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
}

# Scroll the frozen view down near the new rows and leave the selection where
# the author last left off.
$ws.Activate()
$ws.Range("A364").Select()
$excel.ActiveWindow.ScrollRow = 364
$ws.Range("B382").Select()
